$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the last row's B9/C9 cells (removes the shared string "冲刺后后坐力驱动bug" and the "-" reference there)
$ws.Range("B9:C9").Clear()

# Update the active selection to E14 (as seen in the saved sheetView)
$ws.Range("E14").Select()
